# fix akurasi yg salah
# Correct the column headers ("Name" -> "Nama", "gameplay" -> "GamePlay")
# and the accuracy data used by the three line charts (My Bot 1/2/3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header labels for the three data blocks (B/C/D, F/G/H, J/K/L) ---
$ws.Range("B1").Value = "Nama"
$ws.Range("C1").Value = "GamePlay"

$ws.Range("F1").Value = "Nama"
$ws.Range("G1").Value = "GamePlay"

$ws.Range("J1").Value = "Nama"
$ws.Range("K1").Value = "GamePlay"

# --- Corrected accuracy values (reward) for "My Bot 1" (column D) ---
$botOneValues = @(404.72, 1439.74, 1728.6, 2017.33, 1088.4000000000001, 1223.73, 1109.2, 935.17, 400.62, 326.02999999999997)
for ($i = 0; $i -lt $botOneValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 4).Value = $botOneValues[$i]
}

# --- Corrected accuracy values (reward) for "My Bot 2" (column H) ---
$botTwoValues = @(1163.02, 1474.4, 1594, 1059, 1489.48, 1464.09, 1488, 1582.45, 1391, 1396.81)
for ($i = 0; $i -lt $botTwoValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 8).Value = $botTwoValues[$i]
}

# --- Corrected accuracy values (reward) for "My Bot 3" (column L) ---
$botThreeValues = @(2007.96, 856.26, 337.97, 820.56, 1125.48, 1106.28, 1055.42, 1168.5, 1907.24, 1716.02)
for ($i = 0; $i -lt $botThreeValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 12).Value = $botThreeValues[$i]
}
